{"js": "// Minor update to description:\n//   \"There are three main classes in this application:\"\n// becomes two lines (separated by a manual line break) inside the same\n// paragraph:\n//   \"This is a C# console application.\"\n//   \"There are three main classes:\"\n\nconst body = context.document.body;\n\n// Locate the target paragraph by its original text instead of relying on a\n// fixed paragraph index, so the script is resilient to unrelated changes\n// elsewhere in the document.\nconst originalSentence = \"There are three main classes in this application:\";\nconst matches = body.search(originalSentence, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not find the paragraph to update.\");\n}\n\nconst target = matches.items[0].paragraphs.getFirst();\ntarget.load(\"text\");\nawait context.sync();\n\n// 1) Insert a manual line break right before the existing sentence.\nconst startRange = target.getRange(Word.RangeLocation.start);\nstartRange.insertBreak(Word.BreakType.line, Word.InsertLocation.before);\nawait context.sync();\n\n// 2) Insert the new introductory sentence before that line break.\nconst introRange = target.getRange(Word.RangeLocation.start);\nintroRange.insertText(\"This is a C# console application.\", Word.InsertLocation.before);\nawait context.sync();\n\n// 3) Trim the original sentence down to \"There are three main classes:\"\n//    by removing \" in this application\".\ntarget.load(\"text\");\nawait context.sync();\n\nconst trimMatches = target.search(\" in this application:\", { matchCase: true });\ntrimMatches.load(\"items\");\nawait context.sync();\n\nif (trimMatches.items.length > 0) {\n  trimMatches.items[0].insertText(\":\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Minor update to description:\n#   \"There are three main classes in this application:\"\n# becomes two lines (separated by a manual line break) inside the same\n# paragraph:\n#   \"This is a C# console application.\"\n#   \"There are three main classes:\"\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph by its original text instead of relying on a\n# fixed paragraph index, so the script is resilient to unrelated changes\n# elsewhere in the document.\n$originalSentence = \"There are three main classes in this application:\"\n$target = $null\nforeach ($para in $d.Paragraphs) {\n    $paraText = $para.Range.Text.TrimEnd([char]13, [char]7)\n    if ($paraText -eq $originalSentence) {\n        $target = $para\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the paragraph to update.\"\n}\n\n$start = $target.Range.Start\n$introSentence = \"This is a C# console application.\"\n\n# 1) Insert a manual line break (wdLineBreak = 6) at the start of the\n#    paragraph - this lands as its own <w:r><w:br/></w:r> run.\n$breakRange = $d.Range($start, $start)\n$breakRange.InsertBreak(6)\n\n# 2) Insert the new introductory sentence before that line break.\n$introRange = $d.Range($start, $start)\n$introRange.InsertBefore($introSentence)\n\n# 3) Trim the original sentence down to \"There are three main classes:\" by\n#    removing \" in this application\", scoped to this paragraph only.\n$find = $target.Range.Find\n$find.ClearFormatting()\n$find.Text = \" in this application:\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \":\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
